$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two data rows (old rows 9 and 10); this shifts data up
# and updates the sheet dimension from A1:G10 to A1:G8 automatically.
$ws.Rows("9:10").Delete()

# Row 2
$ws.Range("A2").Value = 100
$ws.Range("B2").Value = "Stephens PLC"
$ws.Range("C2").Value = "Hardware"
$ws.Range("D2").Value = "Laptop"
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 3000
$ws.Range("G2").Value = 18000

# Row 3
$ws.Range("A3").Value = 100
$ws.Range("B3").Value = "Stephens PLC"
$ws.Range("C3").Value = "Hardware"
$ws.Range("D3").Value = "Laptop"
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 3000
$ws.Range("G3").Value = 21000

# Row 4
$ws.Range("A4").Value = 100
$ws.Range("B4").Value = "Stephens PLC"
$ws.Range("C4").Value = "Software"
$ws.Range("D4").Value = "Office 365"
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = 200
$ws.Range("G4").Value = 1800

# Row 5
$ws.Range("A5").Value = 100
$ws.Range("B5").Value = "Stephens PLC"
$ws.Range("C5").Value = "Software"
$ws.Range("D5").Value = "Service pack"
$ws.Range("E5").Value = 14
$ws.Range("F5").Value = 100
$ws.Range("G5").Value = 1400

# Row 6
$ws.Range("A6").Value = 100
$ws.Range("B6").Value = "Stephens PLC"
$ws.Range("C6").Value = "Hardware"
$ws.Range("D6").Value = "Router"
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 200
$ws.Range("G6").Value = 2000

# Row 7
$ws.Range("A7").Value = 100
$ws.Range("B7").Value = "Stephens PLC"
$ws.Range("C7").Value = "Hardware"
$ws.Range("D7").Value = "Desktop"
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 2488
$ws.Range("G7").Value = 24880

# Row 8 - blank detail columns, keep only the grand total
$ws.Range("A8").Value = ""
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("G8").Value = 69080
